{"js": "// Map of old cell text -> new cell text for the division-table update.\nconst replacements = [\n  [\"340\u00f74=85, 0\", \"898\u00f78=112, 2\"],\n  [\"863\u00f76=143, 5\", \"121\u00f75=24, 1\"],\n  [\"226\u00f72=113, 0\", \"910\u00f74=227, 2\"],\n  [\"934\u00f74=233, 2\", \"887\u00f72=443, 1\"],\n  [\"105\u00f72=52, 1\", \"407\u00f72=203, 1\"],\n  [\"256\u00f76=42, 4\", \"732\u00f79=81, 3\"],\n  [\"375\u00f79=41, 6\", \"525\u00f72=262, 1\"],\n  [\"220\u00f78=27, 4\", \"336\u00f72=168, 0\"],\n  [\"450\u00f72=225, 0\", \"666\u00f72=333, 0\"],\n  [\"822\u00f76=137, 0\", \"554\u00f75=110, 4\"],\n  [\"544\u00f76=90, 4\", \"878\u00f79=97, 5\"],\n  [\"925\u00f72=462, 1\", \"844\u00f75=168, 4\"],\n  [\"979\u00f78=122, 3\", \"181\u00f72=90, 1\"],\n  [\"756\u00f72=378, 0\", \"632\u00f76=105, 2\"],\n  [\"853\u00f79=94, 7\", \"763\u00f75=152, 3\"],\n  [\"531\u00f75=106, 1\", \"322\u00f79=35, 7\"],\n  [\"524\u00f73=174, 2\", \"202\u00f74=50, 2\"],\n  [\"395\u00f77=56, 3\", \"680\u00f79=75, 5\"],\n  [\"829\u00f72=414, 1\", \"534\u00f74=133, 2\"],\n  [\"352\u00f76=58, 4\", \"668\u00f78=83, 4\"],\n  [\"463\u00f73=154, 1\", \"434\u00f76=72, 2\"],\n  [\"703\u00f78=87, 7\", \"849\u00f78=106, 1\"],\n  [\"411\u00f72=205, 1\", \"972\u00f77=138, 6\"],\n  [\"222\u00f72=111, 0\", \"444\u00f72=222, 0\"],\n  [\"167\u00f79=18, 5\", \"691\u00f72=345, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();", "ps1": "# Update the division-table answers to match the new generated output.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('340\u00f74=85, 0', '898\u00f78=112, 2'),\n    @('863\u00f76=143, 5', '121\u00f75=24, 1'),\n    @('226\u00f72=113, 0', '910\u00f74=227, 2'),\n    @('934\u00f74=233, 2', '887\u00f72=443, 1'),\n    @('105\u00f72=52, 1', '407\u00f72=203, 1'),\n    @('256\u00f76=42, 4', '732\u00f79=81, 3'),\n    @('375\u00f79=41, 6', '525\u00f72=262, 1'),\n    @('220\u00f78=27, 4', '336\u00f72=168, 0'),\n    @('450\u00f72=225, 0', '666\u00f72=333, 0'),\n    @('822\u00f76=137, 0', '554\u00f75=110, 4'),\n    @('544\u00f76=90, 4', '878\u00f79=97, 5'),\n    @('925\u00f72=462, 1', '844\u00f75=168, 4'),\n    @('979\u00f78=122, 3', '181\u00f72=90, 1'),\n    @('756\u00f72=378, 0', '632\u00f76=105, 2'),\n    @('853\u00f79=94, 7', '763\u00f75=152, 3'),\n    @('531\u00f75=106, 1', '322\u00f79=35, 7'),\n    @('524\u00f73=174, 2', '202\u00f74=50, 2'),\n    @('395\u00f77=56, 3', '680\u00f79=75, 5'),\n    @('829\u00f72=414, 1', '534\u00f74=133, 2'),\n    @('352\u00f76=58, 4', '668\u00f78=83, 4'),\n    @('463\u00f73=154, 1', '434\u00f76=72, 2'),\n    @('703\u00f78=87, 7', '849\u00f78=106, 1'),\n    @('411\u00f72=205, 1', '972\u00f77=138, 6'),\n    @('222\u00f72=111, 0', '444\u00f72=222, 0'),\n    @('167\u00f79=18, 5', '691\u00f72=345, 1'),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # 1 = wdFindContinue, 2 = wdReplaceAll\n    $ok = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $ok) {\n        throw \"Text not found: $oldText\"\n    }\n}"}
